$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 239, shifting the existing rows 239-336 down to 241-338.
$ws.Rows.Item(239).Insert()
$ws.Rows.Item(239).Insert()

# Populate the new row 239 with its data.
$ws.Range("A239").Value = 5
$ws.Range("B239").Value = "Macroferia Regional de Talca"
$ws.Range("C239").Value = "Maule"
$ws.Range("D239").Value = 44784
$ws.Range("E239").Value = 7
$ws.Range("F239").Value = 100112006
$ws.Range("G239").Value = "Repollo"
$ws.Range("H239").Value = "Crespo record"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 2000
$ws.Range("K239").Value = 1200
$ws.Range("L239").Value = 1200
$ws.Range("M239").Value = 1200
$ws.Range("N239").Value = "$/unidad"
$ws.Range("O239").Value = "Región del Maule"
$ws.Range("P239").Value = 1200
$ws.Range("Q239").Value = 1
$ws.Range("R239").Value = "Hortaliza"

# Populate the new row 240 with its data.
$ws.Range("A240").Value = 5
$ws.Range("B240").Value = "Macroferia Regional de Talca"
$ws.Range("C240").Value = "Maule"
$ws.Range("D240").Value = 44784
$ws.Range("E240").Value = 7
$ws.Range("F240").Value = 100112006
$ws.Range("G240").Value = "Repollo"
$ws.Range("H240").Value = "Crespo record"
$ws.Range("I240").Value = "Segunda"
$ws.Range("J240").Value = 3000
$ws.Range("K240").Value = 1000
$ws.Range("L240").Value = 1000
$ws.Range("M240").Value = 1000
$ws.Range("N240").Value = "$/unidad"
$ws.Range("O240").Value = "Región del Maule"
$ws.Range("P240").Value = 1000
$ws.Range("Q240").Value = 1
$ws.Range("R240").Value = "Hortaliza"
